$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 11: D11 becomes a formula 45000+210000 ---
$ws.Range("D11").Formula = "=45000+210000"

# --- Row 12: C12 formula gains +10009000 term ---
$ws.Range("C12").Formula = "=14625000+360000+28500000+8800000+10009000"

# --- Row 14: new entry "SALES - cash/retail" ---
$ws.Range("B14").Value = "SALES - cash/retail"
$ws.Range("C14").Formula = "=41859025-27669525-10009000"

# --- Row 15: new entry "SETOR KE BANK" -> actually "JASON - school fee" ---
$ws.Range("B15").Value = "JASON - school fee"
$ws.Range("D15").Value = 11000000

# --- Row 16: new entry "SELISIH - lebih" ---
$ws.Range("B16").Value = "SELISIH - lebih"
$ws.Range("C16").Value = 12000

# --- Row 17: new entry "SETOR KE BANK" ---
$ws.Range("B17").Value = "SETOR KE BANK"
$ws.Range("D17").Formula = "=30000000"

# --- Row 18: new date row, Wages Expense ---
$ws.Range("A18").Value = 44314
$ws.Range("B18").Value = "Wages Expense"
$ws.Range("D18").Formula = "=45000+180000"

# --- Row 19: new entry "BELI abon" ---
$ws.Range("B19").Value = "BELI abon"
$ws.Range("D19").Value = 75000

# --- Row 20: new entry "PAKET LEBARAN" ---
$ws.Range("B20").Value = "PAKET LEBARAN"
$ws.Range("D20").Formula = "=3500000+51000"

# --- Row 21: new entry "TRANSFER BCA" ---
$ws.Range("B21").Value = "TRANSFER BCA"
$ws.Range("D21").Formula = "=100000+5455530+1260000+1700000+4200000+300000+583000"

# --- Row 22: new entry "BONUS FEE CANVASER 2020" ---
$ws.Range("B22").Value = "BONUS FEE CANVASER 2020"
$ws.Range("D22").Formula = "=47787970"

# --- Row 23: new entry "A/R" ---
$ws.Range("B23").Value = "A/R"
$ws.Range("C23").Formula = "=428500+22875000+2175000+2790000+9560000+12800000+750000+940000+925000+10903000"

# --- Row 24: new entry "SALES - cash/retail" ---
$ws.Range("B24").Value = "SALES - cash/retail"
$ws.Range("C24").Formula = "=11152975+9028025-10903000"

# --- Row 25: new entry "SETOR KE BANK" ---
$ws.Range("B25").Value = "SETOR KE BANK"
$ws.Range("D25").Value = 8000000

# --- Row 26: new date row, Wages Expense ---
$ws.Range("A26").Value = 44315
$ws.Range("B26").Value = "Wages Expense"
$ws.Range("D26").Formula = "=45000"

# --- Update the saved view: frozen pane top-left cell and active selection ---
$ws.Activate()
$ws.Range("C25").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
